$wb = $excel.ActiveWorkbook

# --- Master Penjualan sheet: change buyer name "Adib" -> "Cahya" for row 2 ---
$wsMaster = $wb.Worksheets.Item("Master Penjualan")
$wsMaster.Range("B2").Value = "Cahya"

# Update the selected/active cell to reflect the workbook state (C9 -> C10)
$wsMaster.Range("C10").Select()

# --- Update window position (bookViews: xWindow/yWindow) ---
$win = $wb.Windows.Item(1)
$win.Left = 3996
$win.Top = 2148
